$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Edn3"
$ws.Range("C2").Value = "Ednrb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.012975
$ws.Range("H2").Value = 0.038925
$ws.Range("I2").Value = 0.004878150260562778
$ws.Range("J2").Value = 0.004878150260562778
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 70.24576333333334
$ws.Range("N2").Value = 210.73729
$ws.Range("O2").Value = 0.8038202726758664
$ws.Range("P2").Value = 0.8038202726758664
$ws.Range("Q2").Value = 0.9114387792500002
$ws.Range("R2").Value = 8.202949013250002
$ws.Range("S2").Value = 0.003921156072599422
$ws.Range("T2").Value = 0.003921156072599422

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Edn3"
$ws.Range("C3").Value = "Ednrb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.012975
$ws.Range("H3").Value = 0.038925
$ws.Range("I3").Value = 0.004878150260562778
$ws.Range("J3").Value = 0.004878150260562778
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.152389
$ws.Range("N3").Value = 0.457167
$ws.Range("O3").Value = 0.001743782994449666
$ws.Range("P3").Value = 0.001743782994449667
$ws.Range("Q3").Value = 0.001977247275
$ws.Range("R3").Value = 0.017795225475
$ws.Range("S3").Value = 0.000008506435468739582
$ws.Range("T3").Value = 0.000008506435468739582

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Edn3"
$ws.Range("C4").Value = "Ednrb"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.012975
$ws.Range("H4").Value = 0.038925
$ws.Range("I4").Value = 0.004878150260562778
$ws.Range("J4").Value = 0.004878150260562778
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.0875
$ws.Range("N4").Value = 9.2625
$ws.Range("O4").Value = 0.03533017471971957
$ws.Range("P4").Value = 0.03533017471971957
$ws.Range("Q4").Value = 0.0400603125
$ws.Range("R4").Value = 0.3605428125
$ws.Range("S4").Value = 0.0001723459010147285
$ws.Range("T4").Value = 0.0001723459010147285

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Edn3"
$ws.Range("C5").Value = "Ednrb"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.012975
$ws.Range("H5").Value = 0.038925
$ws.Range("I5").Value = 0.004878150260562778
$ws.Range("J5").Value = 0.004878150260562778
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.90423533333333
$ws.Range("N5").Value = 41.712706
$ws.Range("O5").Value = 0.1591057696099644
$ws.Range("P5").Value = 0.1591057696099644
$ws.Range("Q5").Value = 0.18040745345
$ws.Range("R5").Value = 1.62366708105
$ws.Range("S5").Value = 0.000776141851479889
$ws.Range("T5").Value = 0.000776141851479889

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Edn3"
$ws.Range("C6").Value = "Ednrb"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.646844666666667
$ws.Range("H6").Value = 7.940534
$ws.Range("I6").Value = 0.9951218497394373
$ws.Range("J6").Value = 0.9951218497394373
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 70.24576333333334
$ws.Range("N6").Value = 210.73729
$ws.Range("O6").Value = 0.8038202726758664
$ws.Range("P6").Value = 0.8038202726758664
$ws.Range("Q6").Value = 185.9296240347622
$ws.Range("R6").Value = 1673.36661631286
$ws.Range("S6").Value = 0.7998991166032671
$ws.Range("T6").Value = 0.7998991166032671

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Edn3"
$ws.Range("C7").Value = "Ednrb"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.646844666666667
$ws.Range("H7").Value = 7.940534
$ws.Range("I7").Value = 0.9951218497394373
$ws.Range("J7").Value = 0.9951218497394373
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.152389
$ws.Range("N7").Value = 0.457167
$ws.Range("O7").Value = 0.001743782994449666
$ws.Range("P7").Value = 0.001743782994449667
$ws.Range("Q7").Value = 0.4033500119086667
$ws.Range("R7").Value = 3.630150107178
$ws.Range("S7").Value = 0.001735276558980927
$ws.Range("T7").Value = 0.001735276558980927

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Edn3"
$ws.Range("C8").Value = "Ednrb"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.646844666666667
$ws.Range("H8").Value = 7.940534
$ws.Range("I8").Value = 0.9951218497394373
$ws.Range("J8").Value = 0.9951218497394373
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.0875
$ws.Range("N8").Value = 9.2625
$ws.Range("O8").Value = 0.03533017471971957
$ws.Range("P8").Value = 0.03533017471971957
$ws.Range("Q8").Value = 8.172132908333333
$ws.Range("R8").Value = 73.54919617499999
$ws.Range("S8").Value = 0.03515782881870484
$ws.Range("T8").Value = 0.03515782881870484

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Edn3"
$ws.Range("C9").Value = "Ednrb"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.646844666666667
$ws.Range("H9").Value = 7.940534
$ws.Range("I9").Value = 0.9951218497394373
$ws.Range("J9").Value = 0.9951218497394373
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 13.90423533333333
$ws.Range("N9").Value = 41.712706
$ws.Range("O9").Value = 0.1591057696099644
$ws.Range("P9").Value = 0.1591057696099644
$ws.Range("Q9").Value = 36.80235113611155
$ws.Range("R9").Value = 331.221160225004
$ws.Range("S9").Value = 0.1583296277584845
$ws.Range("T9").Value = 0.1583296277584845
